$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# New footstep / land audio rows (ids 1000-1010), appended below the
# existing SFX_1..SFX_3 rows (row 5 previously held the now-removed SFX_4 entry).

# 1) Fill the audioSFXType (column C) values first, top to bottom.
$ws.Cells.Item(5, 3).Value  = "PlayerFootsteps_01"
$ws.Cells.Item(6, 3).Value  = "PlayerFootsteps_02"
$ws.Cells.Item(7, 3).Value  = "PlayerFootsteps_03"
$ws.Cells.Item(8, 3).Value  = "PlayerFootsteps_04"
$ws.Cells.Item(9, 3).Value  = "PlayerFootsteps_05"
$ws.Cells.Item(10, 3).Value = "PlayerFootsteps_06"
$ws.Cells.Item(11, 3).Value = "PlayerFootsteps_07"
$ws.Cells.Item(12, 3).Value = "PlayerFootsteps_08"
$ws.Cells.Item(13, 3).Value = "PlayerFootsteps_09"
$ws.Cells.Item(14, 3).Value = "PlayerFootsteps_010"
$ws.Cells.Item(15, 3).Value = "PlayerLand"

# 2) Fill the name (column B) values: first footstep and the land entry,
#    then backfill the remaining footstep rows in order.
$ws.Cells.Item(5, 2).Value  = "Player_Footstep_01"
$ws.Cells.Item(15, 2).Value = "Player_Land"
$ws.Cells.Item(6, 2).Value  = "Player_Footstep_02"
$ws.Cells.Item(7, 2).Value  = "Player_Footstep_03"
$ws.Cells.Item(8, 2).Value  = "Player_Footstep_04"
$ws.Cells.Item(9, 2).Value  = "Player_Footstep_05"
$ws.Cells.Item(10, 2).Value = "Player_Footstep_06"
$ws.Cells.Item(11, 2).Value = "Player_Footstep_07"
$ws.Cells.Item(12, 2).Value = "Player_Footstep_08"
$ws.Cells.Item(13, 2).Value = "Player_Footstep_09"
$ws.Cells.Item(14, 2).Value = "Player_Footstep_10"

# 3) Fill the id (column A) numeric values.
$ws.Cells.Item(5, 1).Value  = 1000
$ws.Cells.Item(6, 1).Value  = 1001
$ws.Cells.Item(7, 1).Value  = 1002
$ws.Cells.Item(8, 1).Value  = 1003
$ws.Cells.Item(9, 1).Value  = 1004
$ws.Cells.Item(10, 1).Value = 1005
$ws.Cells.Item(11, 1).Value = 1006
$ws.Cells.Item(12, 1).Value = 1007
$ws.Cells.Item(13, 1).Value = 1008
$ws.Cells.Item(14, 1).Value = 1009
$ws.Cells.Item(15, 1).Value = 1010

# Column B is noticeably wider than default to fit the new longer names.
# (20.5 is the input that this engine's pixel-rounded ColumnWidth storage
# resolves to the closest achievable value to the authored 21.33203125.)
$ws.Columns.Item(2).ColumnWidth = 20.5

# Leave the active selection where the author last left off.
$ws.Range("I13").Select()
